# "Correction in the number swap file"
#
# The individual algorithm sheets (BubbleSortClassico, BubbleSortMelhorado,
# InsertionSort) had stale/incorrect swap-count totals for the three largest
# sample sizes (100000, 500000, 1000000) in their "Aleatório" (random) table.
# The "Geral" overview sheet already carried the corrected totals, so this
# fixes the individual sheets to match.
#
# InsertionSort additionally gets two blank rows inserted above its table
# (shifting the whole table down from rows 3-10 to rows 5-12) as part of the
# same edit.

$wb = $excel.ActiveWorkbook

# --- Bubble Sort Classico & Bubble Sort Melhorado: fix the last 3 rows ---
foreach ($name in @("BubbleSortClassico", "BubbleSortMelhorado")) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("D8").Value = 4999950000
    $ws.Range("B9").Value = 62476565906
    $ws.Range("D9").Value = 124999750000
    $ws.Range("B10").Value = 249875026973
    $ws.Range("D10").Value = 499999500000
}

# --- Insertion Sort: insert two rows above the table, then fix last 3 rows ---
$wsIns = $wb.Worksheets.Item("InsertionSort")
$wsIns.Range("A3:D4").EntireRow.Insert()

$wsIns.Range("D10").Value = 4999950000
$wsIns.Range("B11").Value = 62476565906
$wsIns.Range("D11").Value = 124999750000
$wsIns.Range("B12").Value = 249875026973
$wsIns.Range("D12").Value = 499999500000
